# Weekly Ají (pepper) price update for "Terminal La Palmera de La Serena".
# A new weekly observation is inserted as row 175 (shifting the existing
# rows 175:222 down to 176:223), matching the "Fruta / hortaliza, semanal"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175; existing rows shift down one.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A175").Value = 8
$ws.Range("B175").Value = "Terminal La Palmera de La Serena"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44722
$ws.Range("E175").Value = 4
$ws.Range("F175").Value = 100112021
$ws.Range("G175").Value = "Ají"
$ws.Range("H175").Value = "Inferno"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 600
$ws.Range("K175").Value = 21000
$ws.Range("L175").Value = 22000
$ws.Range("M175").Value = 21500
$ws.Range("N175").Value = "`$/caja 12 kilos"
$ws.Range("O175").Value = "Región de Arica y Parinacota"
$ws.Range("P175").Value = 1792
$ws.Range("Q175").Value = 12
$ws.Range("R175").Value = "Hortaliza"
